$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Drag-and-drop" task (row 10) as DEFERRED in the Status column (D)
$ws.Range("D10").Value = "DEFERRED"

# Update the active cell selection to D16
$ws.Range("D16").Select()
